$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12
$ws.Range("I12").Value = 4188377156.3100057
$ws.Range("J12").Value = 35880294

# Row 13
$ws.Range("I13").Value = 1012006300.0300001
$ws.Range("J13").Value = 956934340.60000002

# Row 14
$ws.Range("I14").Value = -44319159.290000051
$ws.Range("J14").Value = 146268235.09999999

# Row 16
$ws.Range("I16").Value = -162861893.56999999
$ws.Range("J16").Value = -193292161.30000001

# Row 18 - I18 becomes a formula, J18 stays a formula (values recompute)
$ws.Range("I18").Formula = "=SUM(I12:I17)"
$ws.Range("J18").Formula = "=SUM(J12:J17)"

# Row 19
$ws.Range("I19").Value = -1160500000.0000002

# Row 21 - I21 becomes a formula, J21 stays a formula
$ws.Range("I21").Formula = "=SUM(I18:I20)"
$ws.Range("J21").Formula = "=SUM(J18:J20)"

# Row 23 - formulas already present, values will recompute automatically

# Row 25 - formulas already present, values will recompute automatically

# Row 26
$ws.Range("I26").Value = 1010658958.9880759
$ws.Range("J26").Value = 1031977291

# Row 28 - formula already present (J28), recomputed automatically

# Row 29 - formula already present (I29), recomputed automatically

$excel.Calculate()
